$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1486.3529
$ws.Range("I100").Value = 1080.9286
$ws.Range("J100").Value = 3378.3333
$ws.Range("K100").Value = 1080.9286
$ws.Range("L100").Value = 3378.3333
$ws.Range("M100").Value = -539.9286
$ws.Range("N100").Value = -4460.3333
$ws.Range("H112").Value = 3535928.8
$ws.Range("J112").Value = 3721925
$ws.Range("L112").Value = 11165775
$ws.Range("N112").Value = -11167991
$ws.Range("H116").Value = 5791.7646
$ws.Range("I116").Value = 6362.8887
$ws.Range("K116").Value = 6362.8887
$ws.Range("M116").Value = -2920.8887
$ws.Range("H137").Value = 13838.538
$ws.Range("I137").Value = 22738.8
$ws.Range("J137").Value = 8275.875
$ws.Range("K137").Value = 68216.39999999999
$ws.Range("L137").Value = 24827.625
$ws.Range("M137").Value = -65666.39999999999
$ws.Range("N137").Value = -29927.625
$ws.Range("H138").Value = 2139.9292
$ws.Range("J138").Value = 2706.9434
$ws.Range("L138").Value = 8120.8302
$ws.Range("N138").Value = -18400.8302

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13718.68
$ws.Range("I32").Value = 12952.292
$ws.Range("J32").Value = 32112
$ws.Range("K32").Value = 12952.292
$ws.Range("L32").Value = 32112
$ws.Range("M32").Value = -12665.292
$ws.Range("N32").Value = -32686
$ws.Range("H61").Value = 6097.6577
$ws.Range("I61").Value = 5078.7
$ws.Range("K61").Value = 5078.7
$ws.Range("M61").Value = -4866.7
$ws.Range("H74").Value = 4103.393
$ws.Range("I74").Value = 4103.393
$ws.Range("K74").Value = 4103.393
$ws.Range("M74").Value = -3229.393
$ws.Range("H77").Value = 4103.393
$ws.Range("I77").Value = 4103.393
$ws.Range("K77").Value = 20516.965
$ws.Range("M77").Value = -16148.965
$ws.Range("H132").Value = 4419.442
$ws.Range("I132").Value = 2875.2307
$ws.Range("K132").Value = 8625.6921
$ws.Range("M132").Value = -6095.6921
$ws.Range("H136").Value = 6097.6577
$ws.Range("I136").Value = 5078.7
$ws.Range("K136").Value = 15236.1
$ws.Range("M136").Value = -12686.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1837.4615
$ws.Range("J86").Value = 2667.2
$ws.Range("L86").Value = 2667.2
$ws.Range("N86").Value = -4913.2
$ws.Range("H89").Value = 1837.4615
$ws.Range("J89").Value = 2667.2
$ws.Range("L89").Value = 13336
$ws.Range("N89").Value = -24568
$ws.Range("H107").Value = 1316.8125
$ws.Range("J107").Value = 1418.4166
$ws.Range("L107").Value = 1418.4166
$ws.Range("N107").Value = -5258.4166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1874.9
$ws.Range("I16").Value = 1919.1666
$ws.Range("J16").Value = 1808.5
$ws.Range("K16").Value = 1919.1666
$ws.Range("L16").Value = 1808.5
$ws.Range("M16").Value = -1632.1666
$ws.Range("N16").Value = -2382.5
$ws.Range("H31").Value = 56800.215
$ws.Range("I31").Value = 3902
$ws.Range("K31").Value = 3902
$ws.Range("M31").Value = -3607
$ws.Range("H34").Value = 56800.215
$ws.Range("I34").Value = 3902
$ws.Range("K34").Value = 3902
$ws.Range("M34").Value = -3700
$ws.Range("H58").Value = 3646.8
$ws.Range("I58").Value = 3740.3333
$ws.Range("K58").Value = 3740.3333
$ws.Range("M58").Value = -3537.3333
$ws.Range("H68").Value = 257500
$ws.Range("J68").Value = 257500
$ws.Range("L68").Value = 257500
$ws.Range("N68").Value = -258998
$ws.Range("H71").Value = 257500
$ws.Range("J71").Value = 257500
$ws.Range("L71").Value = 772500
$ws.Range("N71").Value = -779988
$ws.Range("H88").Value = 11936.25
$ws.Range("J88").Value = 16995
$ws.Range("L88").Value = 16995
$ws.Range("N88").Value = -17807
$ws.Range("H91").Value = 11936.25
$ws.Range("J91").Value = 16995
$ws.Range("L91").Value = 16995
$ws.Range("N91").Value = -19803
$ws.Range("H99").Value = 2747.1333
$ws.Range("I99").Value = 2732.2
$ws.Range("K99").Value = 2732.2
$ws.Range("M99").Value = -1234.2
$ws.Range("H107").Value = 3671.5217
$ws.Range("I107").Value = 429.42105
$ws.Range("J107").Value = 5953
$ws.Range("K107").Value = 429.42105
$ws.Range("L107").Value = 5953
$ws.Range("M107").Value = 1490.57895
$ws.Range("N107").Value = -9793
$ws.Range("H113").Value = 1874.9
$ws.Range("I113").Value = 1919.1666
$ws.Range("J113").Value = 1808.5
$ws.Range("K113").Value = 1919.1666
$ws.Range("L113").Value = 1808.5
$ws.Range("M113").Value = 250.8334
$ws.Range("N113").Value = -6148.5
$ws.Range("H126").Value = 2747.1333
$ws.Range("I126").Value = 2732.2
$ws.Range("K126").Value = 8196.599999999999
$ws.Range("M126").Value = -5726.599999999999
$ws.Range("H132").Value = 3148.5334
$ws.Range("I132").Value = 2979.8845
$ws.Range("J132").Value = 4244.75
$ws.Range("K132").Value = 8939.6535
$ws.Range("L132").Value = 12734.25
$ws.Range("M132").Value = -6409.6535
$ws.Range("N132").Value = -17794.25
$ws.Range("H136").Value = 3646.8
$ws.Range("I136").Value = 3740.3333
$ws.Range("K136").Value = 11220.9999
$ws.Range("M136").Value = -8670.999899999999
$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360
$ws.Range("H141").Value = 128295.97
$ws.Range("J141").Value = 128295.97
$ws.Range("L141").Value = 128295.97
$ws.Range("N141").Value = -138655.97

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 2033.3334
$ws.Range("J45").Value = 2500
$ws.Range("L45").Value = 7500
$ws.Range("N45").Value = -8564
$ws.Range("H132").Value = 5893.2964
$ws.Range("I132").Value = 6680.6816
$ws.Range("J132").Value = 2428.8
$ws.Range("K132").Value = 60126.1344
$ws.Range("L132").Value = 21859.2
$ws.Range("M132").Value = -57596.1344
$ws.Range("N132").Value = -26919.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 26134.768
$ws.Range("I132").Value = 29436.139
$ws.Range("J132").Value = 9156.286
$ws.Range("K132").Value = 88308.417
$ws.Range("L132").Value = 27468.858
$ws.Range("M132").Value = -85778.417
$ws.Range("N132").Value = -32528.858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7176.381
$ws.Range("I132").Value = 3325.7
$ws.Range("K132").Value = 9977.099999999999
$ws.Range("M132").Value = -7447.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 32397.6
$ws.Range("J52").Value = 23994.5
$ws.Range("L52").Value = 23994.5
$ws.Range("N52").Value = -24446.5
$ws.Range("H69").Value = 33914
$ws.Range("J69").Value = 33914
$ws.Range("L69").Value = 33914
$ws.Range("N69").Value = -35412
$ws.Range("H72").Value = 33914
$ws.Range("J72").Value = 33914
$ws.Range("L72").Value = 101742
$ws.Range("N72").Value = -109230
$ws.Range("H132").Value = 2179.2222
$ws.Range("I132").Value = 951.75
$ws.Range("K132").Value = 2855.25
$ws.Range("M132").Value = -325.25
